$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '50.999.89'
$ws.Range('E2').Value = '  -1.61%  '

$ws.Range('D3').Value = '2.935.88'
$ws.Range('E3').Value = '  -1.86%  '

$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').Style = "Normal"
$ws.Range('E4').Value = '  +0.09%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '378.19'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -0.10%  '

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '102.26'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -2.97%  '

$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.538'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  -1.54%  '

$ws.Range('E8').Value = '  +0.01%  '

$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.586'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  -2.68%  '

$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '36.51'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  -3.20%  '

$ws.Range('E12').Value = '  -1.22%  '

$ws.Range('D13').Value = '3.405.11'
$ws.Range('E13').Value = '  -1.51%  '

$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '17.99'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -4.17%  '

$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '7.37'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  -2.27%  '

$ws.Range('D16').Value = '2.930.23'
$ws.Range('E16').Value = '  -2.43%  '

$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.982'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  +1.78%  '

$ws.Range('D18').Value = '50.992.28'
$ws.Range('E18').Value = '  -1.81%  '

$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '3.19'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  -8.34%  '

$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '7.13'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -4.41%  '

$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '12.53'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -5.37%  '

$ws.Range('D22').Value = '0.0₃0951'
$ws.Range('E22').Value = '  -1.20%  '

$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '68.48'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -0.62%  '

$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '261.66'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  -0.99%  '

$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.90'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +3.66%  '

$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '8.26'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  +9.42%  '

$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '7.69'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +2.99%  '

$ws.Range('E28').Value = '  -1.74%  '

$ws.Range('B29').Value = 'Hedera'
$ws.Range('C29').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '0.113'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  +7.62%  '

$ws.Range('B30').Value = 'Dai'
$ws.Range('C30').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '1.00'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -0.04%  '

$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '25.61'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -2.35%  '

$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '9.79'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  -1.95%  '

$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '34.20'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  -1.80%  '

$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.0455'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  +3.58%  '

$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '50.33'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  -2.23%  '

$ws.Range('E36').Value = '  -1.45%  '

$ws.Range('E37').Value = '  +0.01%  '

$ws.Range('E38').Value = '  -4.02%  '

$ws.Range('E39').Value = '  -4.93%  '

$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '16.77'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  -4.19%  '

$ws.Range('E41').Value = '  -0.90%  '

$ws.Range('E42').Value = '  -5.04%  '

$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '120.80'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -2.72%  '

$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '21.26'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -5.11%  '

$ws.Range('E45').Value = '  -1.62%  '

$ws.Range('E46').Value = '  +2.47%  '

$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.272'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -4.27%  '

$ws.Range('B48').Value = 'Maker'
$ws.Range('C48').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D48').Value = '2.005.43'
$ws.Range('E48').Value = '  -1.97%  '

$ws.Range('B49').Value = 'NEARProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '3.22'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -1.24%  '

$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.0346'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  +2.14%  '

$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.478'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +12.09%  '

